$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) rows 2 & 3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-20 01:18:31"
$wsOverview.Range("G3").Value = "2016-10-20 01:18:31"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-20 01:18:20"
$wsZhCn.Range("H3").Value = "2016-10-20 01:18:20"
$wsZhCn.Range("K2").Value = "2016-10-20 01:19:00"
$wsZhCn.Range("K3").Value = "2016-10-20 01:19:00"

# de-de sheet: Priority (E), Correspond Handback DateTime (K)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("K2").Value = "2016-10-20 01:19:17"
$wsDeDe.Range("K3").Value = "2016-10-20 01:19:17"
